# Replace the multiplication-fact answers in the table cells with their
# new values, as per the commit's regenerated-output diff.
$d = $word.ActiveDocument

$d.Content.Find.Execute("82×93=7626", $true, $false, $false, $false, $false, $true, 1, $false, "31×89=2759", 2) | Out-Null
$d.Content.Find.Execute("96×17=1632", $true, $false, $false, $false, $false, $true, 1, $false, "35×94=3290", 2) | Out-Null
$d.Content.Find.Execute("65×32=2080", $true, $false, $false, $false, $false, $true, 1, $false, "81×17=1377", 2) | Out-Null
$d.Content.Find.Execute("82×63=5166", $true, $false, $false, $false, $false, $true, 1, $false, "55×16=880", 2) | Out-Null
$d.Content.Find.Execute("12×60=720", $true, $false, $false, $false, $false, $true, 1, $false, "84×77=6468", 2) | Out-Null
$d.Content.Find.Execute("11×18=198", $true, $false, $false, $false, $false, $true, 1, $false, "23×79=1817", 2) | Out-Null
$d.Content.Find.Execute("39×86=3354", $true, $false, $false, $false, $false, $true, 1, $false, "48×26=1248", 2) | Out-Null
$d.Content.Find.Execute("57×39=2223", $true, $false, $false, $false, $false, $true, 1, $false, "15×96=1440", 2) | Out-Null
$d.Content.Find.Execute("55×13=715", $true, $false, $false, $false, $false, $true, 1, $false, "66×59=3894", 2) | Out-Null
$d.Content.Find.Execute("23×52=1196", $true, $false, $false, $false, $false, $true, 1, $false, "17×39=663", 2) | Out-Null
$d.Content.Find.Execute("98×13=1274", $true, $false, $false, $false, $false, $true, 1, $false, "97×37=3589", 2) | Out-Null
$d.Content.Find.Execute("56×36=2016", $true, $false, $false, $false, $false, $true, 1, $false, "99×29=2871", 2) | Out-Null
$d.Content.Find.Execute("49×64=3136", $true, $false, $false, $false, $false, $true, 1, $false, "84×16=1344", 2) | Out-Null
$d.Content.Find.Execute("38×70=2660", $true, $false, $false, $false, $false, $true, 1, $false, "85×79=6715", 2) | Out-Null
$d.Content.Find.Execute("34×62=2108", $true, $false, $false, $false, $false, $true, 1, $false, "59×81=4779", 2) | Out-Null
$d.Content.Find.Execute("82×77=6314", $true, $false, $false, $false, $false, $true, 1, $false, "36×74=2664", 2) | Out-Null
$d.Content.Find.Execute("62×81=5022", $true, $false, $false, $false, $false, $true, 1, $false, "40×33=1320", 2) | Out-Null
$d.Content.Find.Execute("42×99=4158", $true, $false, $false, $false, $false, $true, 1, $false, "67×19=1273", 2) | Out-Null
$d.Content.Find.Execute("81×40=3240", $true, $false, $false, $false, $false, $true, 1, $false, "72×36=2592", 2) | Out-Null
$d.Content.Find.Execute("25×51=1275", $true, $false, $false, $false, $false, $true, 1, $false, "90×82=7380", 2) | Out-Null
$d.Content.Find.Execute("59×51=3009", $true, $false, $false, $false, $false, $true, 1, $false, "86×31=2666", 2) | Out-Null
$d.Content.Find.Execute("91×21=1911", $true, $false, $false, $false, $false, $true, 1, $false, "19×15=285", 2) | Out-Null
$d.Content.Find.Execute("17×66=1122", $true, $false, $false, $false, $false, $true, 1, $false, "89×38=3382", 2) | Out-Null
$d.Content.Find.Execute("93×88=8184", $true, $false, $false, $false, $false, $true, 1, $false, "96×95=9120", 2) | Out-Null
$d.Content.Find.Execute("98×69=6762", $true, $false, $false, $false, $false, $true, 1, $false, "13×23=299", 2) | Out-Null
